# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 100 (pushing the existing
# rows 100-128 down to 101-129) for "Macroferia Regional de Talca" /
# Arandano (blue).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 100..128 down to 101..129, leaving a blank row 100.
$ws.Rows.Item(100).Insert()

# Populate the newly inserted row 100 with the new record.
$ws.Range("A100").Value = 5
$ws.Range("B100").Value = "Macroferia Regional de Talca"
$ws.Range("C100").Value = "Maule"
$ws.Range("D100").Value = 44985
$ws.Range("E100").Value = 7
$ws.Range("F100").Value = "Fruta"
$ws.Range("G100").Value = 100101
$ws.Range("H100").Value = "Berries"
$ws.Range("I100").Value = 100101001
$ws.Range("J100").Value = "Arándano (blue)"
$ws.Range("K100").Value = "Sin especificar"
$ws.Range("L100").Value = "Primera"
$ws.Range("M100").Value = 120
$ws.Range("N100").Value = 3000
$ws.Range("O100").Value = 3000
$ws.Range("P100").Value = 3000
$ws.Range("Q100").Value = "$/bandeja 2 kilos"
$ws.Range("R100").Value = "Provincia de Curicó"
$ws.Range("S100").Value = 1500
$ws.Range("T100").Value = 2
